$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

# Bulk-update column C (dynamic wait flag) from "NO" to "YES" for rows 6-69
for ($r = 6; $r -le 69; $r++) {
    $ws.Range("C$r").Value = "YES"
}

# Reflect the user's navigation/selection state: Registration tab becomes
# the active sheet, with C2:C69 selected (active cell C2)
$ws.Activate()
$ws.Range("C2:C69").Select()
